$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("总计")
$ws2 = $wb.Worksheets.Item("2022-Q2")

# --- 1. "总计" sheet: insert a new summary row for 2022-Q4 above the
#        existing 2022-Q2 row, then fill in the new totals. ---
$ws1.Rows.Item(2).Insert()

# The inserted row inherited row-1's formatting; reset it then restore the
# index-column style (s="2") by copying it from the row that shifted down.
$ws1.Range("A2:D2").ClearFormats()
$ws1.Range("A3").Copy()
$ws1.Range("A2").PasteSpecial(-4122)

$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q4"
$ws1.Range("C2").Value = 5
$ws1.Range("D2").Value = 0.3

# The old row's index counter advances now that it is the 2nd data row.
$ws1.Range("A3").Value = 1

# --- 2. Add the new "2022-Q4" detail sheet, positioned between "总计"
#        and "2022-Q2" (matches the final tab order). ---
$newSheet = $wb.Worksheets.Add($ws2)
$newSheet.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

$rows = @(
    @("003956", "南方产业智选股票",            "3.60", "85.80",  "4.41", "0.1588", 9),
    @("013028", "恒越品质生活混合",              "1.25", "90.58",  "4.59", "0.0574", 2),
    @("012182", "广发沪港深精选混合A",           "0.84", "87.00",  "4.54", "0.0381", 9),
    @("005646", "中海沪港深多策略灵活配置混合",   "0.79", "105.93", "4.56", "0.0360", 7),
    @("012183", "广发沪港深精选混合C",           "0.19", "87.00",  "4.54", "0.0086", 9)
)

# Columns B-G are stored as plain text in the source data (leading zeros
# in fund codes, fixed-decimal strings) - force text formatting on the
# whole block before assignment so Excel doesn't coerce them to numbers,
# then drop the number-format override again so no stray style sticks to
# the cells once the values are in place.
$textBlock = $newSheet.Range("B2:G6")
$textBlock.NumberFormat = "@"

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = 2 + $r
    $rowData = $rows[$r]

    $newSheet.Cells.Item($row, 1).Value = $r

    for ($c = 0; $c -lt 6; $c++) {
        $newSheet.Cells.Item($row, 2 + $c).Value = $rowData[$c]
    }

    # Last column (仓位排名) is a genuine number.
    $newSheet.Cells.Item($row, 8).Value = $rowData[6]
}

$textBlock.ClearFormats()

# Apply the header style (s="2", bold + border) by copying it from the
# already-styled header cell on "总计".
$ws1.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Apply the index-column style (s="2") to the new sheet's A column too.
$ws1.Range("A2").Copy()
$newSheet.Range("A2:A6").PasteSpecial(-4122)

$excel.CutCopyMode = $false
